$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix B89: was stored as text "4", should become a real number 4 ---
$ws.Cells.Item(89, 2).Value = 4

# --- Append new row 90 with Ying Tang's annotation ---
$ws.Cells.Item(90, 1).Value = "Ying Tang"
# politeness_score for this row stays textual ("3"), matching the source data
$ws.Cells.Item(90, 2).Value = "'3"
$ws.Cells.Item(90, 3).Value = "we do not believe"
$ws.Cells.Item(90, 4).Value = "DIS"
$ws.Cells.Item(90, 5).Value = "MET"
$ws.Cells.Item(90, 6).Value = "216e3c96-70ff-4d1d-bc9b-ae161e0068a3"
$ws.Cells.Item(90, 7).Value = "BJInEZsTb_annotated.xlsx"
$ws.Cells.Item(90, 8).Value = "While this is true, we do not believe is necessarily constitutes a disadvantage of our networks, especially when considering ease of training and reproducibility."
